$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.711.04'
$ws.Range("E2").Value = '  +0.51%  '

$ws.Range("D3").Value = '1.844.63'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.54'
$ws.Range("E5").Value = '  -0.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4318'
$ws.Range("E7").Value = '  +0.97%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3644'
$ws.Range("E8").Value = '  -0.20%  '

# row 9: OKB
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.98'
$ws.Range("E9").Value = '  -2.01%  '

# row 10: Dogecoin
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07335'
$ws.Range("E10").Value = '  +0.63%  '

# row 11: Polygon
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8778'
$ws.Range("E11").Value = '  -1.92%  '

# row 12: Solana
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.72'
$ws.Range("E12").Value = '  +0.45%  '

# row 13: WrappedEther
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.873.05'
$ws.Range("E13").Value = '  +0.34%  '

# row 14: Polkadot
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.348'
$ws.Range("E14").Value = '  -0.63%  '

# row 15: Chainlink
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.523'
$ws.Range("E15").Value = '  -0.69%  '

# row 16: TRON
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06964'
$ws.Range("E16").Value = '  +1.22%  '

# row 17: BinanceUSD
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("E17").Value = '  +0.13%  '

# row 18: Litecoin
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '79.66'
$ws.Range("E18").Value = '  +1.48%  '

# row 19: ShibaInu
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009010'
$ws.Range("E19").Value = '  +1.22%  '

# row 20: Dai
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.13%  '

# row 21: Avalanche
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.40'
$ws.Range("E21").Value = '  -0.92%  '

# row 22: WrappedBTC
$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").Value = '27.672.26'
$ws.Range("E22").Value = '  +0.33%  '

# row 23: Uniswap
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.977'
$ws.Range("E23").Value = '  -0.17%  '

# row 24: Cosmos
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.30'
$ws.Range("E24").Value = '  -2.15%  '

# row 25: WrappedliquidstakedEther2.0
$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = '2.078.98'
$ws.Range("E25").Value = '  -1.44%  '

# row 26: Toncoin
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.992'
$ws.Range("E26").Value = '  -2.44%  '

# row 27: Monero
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.74'
$ws.Range("E27").Value = '  +0.22%  '

# row 28: EthereumClassic
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.59'
$ws.Range("E28").Value = '  +0.89%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '119.49'
$ws.Range("E29").Value = '  +3.08%  '

# row 30: InternetComputer(DFINITY)
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.221'
$ws.Range("E30").Value = '  +0.01%  '

# row 31: LidoDAOToken
$ws.Range("B31").Value = 'LidoDAOToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.870'
$ws.Range("E31").Value = '  +1.96%  '

# row 32: Stellar
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08893'
$ws.Range("E32").Value = '  -0.20%  '

# row 33: ImmutableX
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7534'
$ws.Range("E33").Value = '  -3.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.964'
$ws.Range("E34").Value = '  +0.49%  '

# row 35: Filecoin
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.526'
$ws.Range("E35").Value = '  -0.93%  '

# row 36: ARBITRUM
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.118'
$ws.Range("E36").Value = '  +0.86%  '

# row 37: Frax
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.003'
$ws.Range("E37").Value = '  +0.23%  '

# row 38: Hedera
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05436'
$ws.Range("E38").Value = '  -0.07%  '

# row 39: TrustWalletToken
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.106'
$ws.Range("E39").Value = '  +0.52%  '

# row 40: VeChain
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01936'
$ws.Range("E40").Value = '  +0.66%  '

# row 41: MXToken
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.829'
$ws.Range("E41").Value = '  +2.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1665'
$ws.Range("E42").Value = '  +0.92%  '

# row 43: TheSandbox
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5077'
$ws.Range("E43").Value = '  +0.16%  '

# row 44: FraxShare
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.617'
$ws.Range("E44").Value = '  -3.33%  '

# row 45: Aptos
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.364'
$ws.Range("E45").Value = '  +1.62%  '

# row 46: Cronos
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06548'
$ws.Range("E46").Value = '  -1.19%  '

# row 47: EnergySwap
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.40'
$ws.Range("E47").Value = '  +0.37%  '

# row 48: Quant
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.53'
$ws.Range("E48").Value = '  -0.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.003'
$ws.Range("E49").Value = '  +0.16%  '

# row 50: Decentraland
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4639'
$ws.Range("E50").Value = '  -1.53%  '

# row 51: NEARProtocol
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.640'
$ws.Range("E51").Value = '  +0.41%  '
